# StatusTableOverview.xlsx update
# - DoorController PlantUML review done: mark status changed, add source link,
#   reached states/transitions, and update comment.
# - Tweak "need potential properties" wording to "may need properties" (rows 4 & 9).
# - Fill in missing "Code Edits" confirmation note for Producer Consumer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Train Door Controllers): plantuml reviewed -> fill in details
$doorControllerPaperUrl = "https://rebeca-lang.org/assets/papers/2020/Towards-a-Verification-Driven-Iterative-Development-of-Cyber-Physical-System.pdf"
$ws.Range("C10").Value = $doorControllerPaperUrl
$ws.Hyperlinks.Add($ws.Range("C10"), $doorControllerPaperUrl) | Out-Null
$ws.Range("C10").Style = "Hyperlink"

# Row 4 (Scientific Lab) and Row 9 (Ticket Service): reword status text
$ws.Range("D4").Value = "OK, may need properties"
$ws.Range("D9").Value = "OK, may need properties"

# Rows 5 & 8 (Commit Problem / Sender Receiver): Code Edits cell format
# realigned to match the Bad status already shown in Check-in-Afra
$ws.Range("E5").Style = "Bad"
$ws.Range("E8").Style = "Bad"

# Row 6 (Producer Consumer): Code Edits column was blank, now confirmed
$ws.Range("E6").Value = "need to be confirmed"
$ws.Range("E6").Style = "Neutral"

# Row 10 continued
$ws.Range("D10").Value = "OK, may need properties"
$ws.Range("D10").Style = "Neutral"
$ws.Range("F10").Value = 471
$ws.Range("F10").Style = "Neutral"
$ws.Range("G10").Value = 537
$ws.Range("G10").Style = "Neutral"
$ws.Range("H10").Value = "Changes done"
$ws.Range("H10").Style = "Neutral"
$ws.Range("I10").Value = "Changed to be more alike our own process for translating diagrams."

# Reflect that the user last looked at/selected cell I9, scrolled right
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("I9").Select()
